$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.370.73'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -4.20%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.569.28'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.65%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.001'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3674'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '49.20'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3383'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.79%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.166'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07620'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.04%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.23'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.73%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.056'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.907'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.565.75'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.03%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001132'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.53'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06748'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.002'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.213'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.5335'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.75%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.50'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.31%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.98'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '22.394.30'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.374'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.906'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '145.32'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.976'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '125.53'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.49%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.747.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.72%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.046'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.64%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.241'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.015'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.41%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.21'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.65%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08451'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02536'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2316'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.03%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.537'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.73%  '
$ws.Range('B41').Value = 'Hedera'
$ws.Range('C41').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.06492'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.48%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.304'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.68'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6343'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.36%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.31'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.38%  '
$ws.Range('E46').Value = '  +0.05%  '
$ws.Range('E47').Value = '  -5.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.748'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.100'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.262'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '124.37'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.74%  '
